$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3663.348
$ws.Cells.Item(64, 9).Value = 3388
$ws.Cells.Item(64, 10).Value = 4654.6
$ws.Cells.Item(64, 11).Value = 3388
$ws.Cells.Item(64, 12).Value = 4654.6
$ws.Cells.Item(64, 13).Value = -3140
$ws.Cells.Item(64, 14).Value = -5150.6
$ws.Cells.Item(67, 8).Value = 3663.348
$ws.Cells.Item(67, 9).Value = 3388
$ws.Cells.Item(67, 10).Value = 4654.6
$ws.Cells.Item(67, 11).Value = 3388
$ws.Cells.Item(67, 12).Value = 4654.6
$ws.Cells.Item(67, 13).Value = -2530
$ws.Cells.Item(67, 14).Value = -6370.6
$ws.Cells.Item(68, 8).Value = 90000
$ws.Cells.Item(68, 10).Value = 90000
$ws.Cells.Item(68, 12).Value = 90000
$ws.Cells.Item(68, 14).Value = -91498
$ws.Cells.Item(69, 8).Value = 3679149.2
$ws.Cells.Item(69, 9).Value = 7354941
$ws.Cells.Item(69, 10).Value = 3357.5
$ws.Cells.Item(69, 11).Value = 22064823
$ws.Cells.Item(69, 12).Value = 10072.5
$ws.Cells.Item(69, 13).Value = -22063949
$ws.Cells.Item(69, 14).Value = -11820.5
$ws.Cells.Item(70, 8).Value = 2269.2307
$ws.Cells.Item(70, 9).Value = 2271.1428
$ws.Cells.Item(70, 10).Value = 2267
$ws.Cells.Item(70, 11).Value = 6813.428400000001
$ws.Cells.Item(70, 12).Value = 6801
$ws.Cells.Item(70, 13).Value = -6543.428400000001
$ws.Cells.Item(70, 14).Value = -7341
$ws.Cells.Item(71, 8).Value = 90000
$ws.Cells.Item(71, 10).Value = 90000
$ws.Cells.Item(71, 12).Value = 270000
$ws.Cells.Item(71, 14).Value = -277488
$ws.Cells.Item(72, 8).Value = 3679149.2
$ws.Cells.Item(72, 9).Value = 7354941
$ws.Cells.Item(72, 10).Value = 3357.5
$ws.Cells.Item(72, 11).Value = 66194469
$ws.Cells.Item(72, 12).Value = 30217.5
$ws.Cells.Item(72, 13).Value = -66190101
$ws.Cells.Item(72, 14).Value = -38953.5
$ws.Cells.Item(73, 8).Value = 2269.2307
$ws.Cells.Item(73, 9).Value = 2271.1428
$ws.Cells.Item(73, 10).Value = 2267
$ws.Cells.Item(73, 11).Value = 6813.428400000001
$ws.Cells.Item(73, 12).Value = 6801
$ws.Cells.Item(73, 13).Value = -5877.428400000001
$ws.Cells.Item(73, 14).Value = -8673
$ws.Cells.Item(74, 8).Value = 2225682
$ws.Cells.Item(74, 9).Value = 3128434.5
$ws.Cells.Item(74, 11).Value = 3128434.5
$ws.Cells.Item(74, 13).Value = -3127498.5
$ws.Cells.Item(76, 8).Value = 55559012
$ws.Cells.Item(76, 9).Value = 62502950
$ws.Cells.Item(76, 11).Value = 62502950
$ws.Cells.Item(76, 13).Value = -62502635
$ws.Cells.Item(77, 8).Value = 2225682
$ws.Cells.Item(77, 9).Value = 3128434.5
$ws.Cells.Item(77, 11).Value = 15642172.5
$ws.Cells.Item(77, 13).Value = -15637492.5
$ws.Cells.Item(79, 8).Value = 55559012
$ws.Cells.Item(79, 9).Value = 62502950
$ws.Cells.Item(79, 11).Value = 62502950
$ws.Cells.Item(79, 13).Value = -62501858
$ws.Cells.Item(82, 8).Value = 1760.2354
$ws.Cells.Item(82, 9).Value = 318.5
$ws.Cells.Item(82, 10).Value = 3041.7778
$ws.Cells.Item(82, 11).Value = 955.5
$ws.Cells.Item(82, 12).Value = 9125.3334
$ws.Cells.Item(82, 13).Value = -549.5
$ws.Cells.Item(82, 14).Value = -9937.3334
$ws.Cells.Item(85, 8).Value = 1760.2354
$ws.Cells.Item(85, 9).Value = 318.5
$ws.Cells.Item(85, 10).Value = 3041.7778
$ws.Cells.Item(85, 11).Value = 955.5
$ws.Cells.Item(85, 12).Value = 9125.3334
$ws.Cells.Item(85, 13).Value = 448.5
$ws.Cells.Item(85, 14).Value = -11933.3334
$ws.Cells.Item(137, 8).Value = 57337.277
$ws.Cells.Item(137, 9).Value = 334498
$ws.Cells.Item(137, 10).Value = 1905.1333
$ws.Cells.Item(137, 11).Value = 1003494
$ws.Cells.Item(137, 12).Value = 5715.3999
$ws.Cells.Item(137, 13).Value = -1000944
$ws.Cells.Item(137, 14).Value = -10815.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2068.925
$ws.Cells.Item(61, 9).Value = 2109.5715
$ws.Cells.Item(61, 10).Value = 1784.4
$ws.Cells.Item(61, 11).Value = 2109.5715
$ws.Cells.Item(61, 12).Value = 1784.4
$ws.Cells.Item(61, 13).Value = -1897.5715
$ws.Cells.Item(61, 14).Value = -2208.4
$ws.Cells.Item(74, 8).Value = 58749.312
$ws.Cells.Item(74, 9).Value = 78649.234
$ws.Cells.Item(74, 10).Value = 1260.6666
$ws.Cells.Item(74, 11).Value = 78649.234
$ws.Cells.Item(74, 12).Value = 1260.6666
$ws.Cells.Item(74, 13).Value = -77775.234
$ws.Cells.Item(74, 14).Value = -3008.6666
$ws.Cells.Item(77, 8).Value = 58749.312
$ws.Cells.Item(77, 9).Value = 78649.234
$ws.Cells.Item(77, 10).Value = 1260.6666
$ws.Cells.Item(77, 11).Value = 393246.17
$ws.Cells.Item(77, 12).Value = 6303.333000000001
$ws.Cells.Item(77, 13).Value = -388878.17
$ws.Cells.Item(77, 14).Value = -15039.333
$ws.Cells.Item(136, 8).Value = 2068.925
$ws.Cells.Item(136, 9).Value = 2109.5715
$ws.Cells.Item(136, 10).Value = 1784.4
$ws.Cells.Item(136, 11).Value = 6328.7145
$ws.Cells.Item(136, 12).Value = 5353.200000000001
$ws.Cells.Item(136, 13).Value = -3778.7145
$ws.Cells.Item(136, 14).Value = -10453.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(132, 8).Value = 250022500
$ws.Cells.Item(132, 10).Value = 250022500
$ws.Cells.Item(132, 12).Value = 250022500
$ws.Cells.Item(132, 14).Value = -250032620
$ws.Cells.Item(134, 8).Value = 35788520
$ws.Cells.Item(134, 9).Value = 62500652
$ws.Cells.Item(134, 10).Value = 172344
$ws.Cells.Item(134, 11).Value = 187501956
$ws.Cells.Item(134, 12).Value = 517032
$ws.Cells.Item(134, 13).Value = -187499421
$ws.Cells.Item(134, 14).Value = -522102

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1080
$ws.Cells.Item(16, 10).Value = 1133.3334
$ws.Cells.Item(16, 12).Value = 1133.3334
$ws.Cells.Item(16, 14).Value = -1707.3334
$ws.Cells.Item(31, 8).Value = 49223.684
$ws.Cells.Item(31, 9).Value = 101592.1
$ws.Cells.Item(31, 10).Value = 5583.3335
$ws.Cells.Item(31, 11).Value = 101592.1
$ws.Cells.Item(31, 12).Value = 5583.3335
$ws.Cells.Item(31, 13).Value = -101297.1
$ws.Cells.Item(31, 14).Value = -6173.3335
$ws.Cells.Item(34, 8).Value = 49223.684
$ws.Cells.Item(34, 9).Value = 101592.1
$ws.Cells.Item(34, 10).Value = 5583.3335
$ws.Cells.Item(34, 11).Value = 101592.1
$ws.Cells.Item(34, 12).Value = 5583.3335
$ws.Cells.Item(34, 13).Value = -101390.1
$ws.Cells.Item(34, 14).Value = -5987.3335
$ws.Cells.Item(58, 8).Value = 2482.5938
$ws.Cells.Item(58, 9).Value = 880.1177
$ws.Cells.Item(58, 10).Value = 4298.7334
$ws.Cells.Item(58, 11).Value = 880.1177
$ws.Cells.Item(58, 12).Value = 4298.7334
$ws.Cells.Item(58, 13).Value = -677.1177
$ws.Cells.Item(58, 14).Value = -4704.7334
$ws.Cells.Item(113, 8).Value = 1080
$ws.Cells.Item(113, 10).Value = 1133.3334
$ws.Cells.Item(113, 12).Value = 1133.3334
$ws.Cells.Item(113, 14).Value = -5473.3334
$ws.Cells.Item(132, 8).Value = 1878.6
$ws.Cells.Item(132, 9).Value = 1386.7059
$ws.Cells.Item(132, 11).Value = 4160.1177
$ws.Cells.Item(132, 13).Value = -1630.1177
$ws.Cells.Item(134, 8).Value = 17858616
$ws.Cells.Item(134, 9).Value = 1443.9375
$ws.Cells.Item(134, 10).Value = 41668180
$ws.Cells.Item(134, 11).Value = 4331.8125
$ws.Cells.Item(134, 12).Value = 125004540
$ws.Cells.Item(134, 13).Value = -1796.8125
$ws.Cells.Item(134, 14).Value = -125009610
$ws.Cells.Item(136, 8).Value = 2482.5938
$ws.Cells.Item(136, 9).Value = 880.1177
$ws.Cells.Item(136, 10).Value = 4298.7334
$ws.Cells.Item(136, 11).Value = 2640.3531
$ws.Cells.Item(136, 12).Value = 12896.2002
$ws.Cells.Item(136, 13).Value = -90.35310000000027
$ws.Cells.Item(136, 14).Value = -17996.2002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 8297.35
$ws.Cells.Item(80, 10).Value = 12078.917
$ws.Cells.Item(80, 12).Value = 12078.917
$ws.Cells.Item(80, 14).Value = -14074.917
$ws.Cells.Item(83, 8).Value = 8297.35
$ws.Cells.Item(83, 10).Value = 12078.917
$ws.Cells.Item(83, 12).Value = 60394.585
$ws.Cells.Item(83, 14).Value = -70378.58499999999
$ws.Cells.Item(102, 8).Value = 32404
$ws.Cells.Item(102, 9).Value = 11322.4
$ws.Cells.Item(102, 11).Value = 11322.4
$ws.Cells.Item(102, 13).Value = -9700.4
$ws.Cells.Item(122, 8).Value = 2461.5386
$ws.Cells.Item(122, 9).Value = 2672.7273
$ws.Cells.Item(122, 10).Value = 1300
$ws.Cells.Item(122, 11).Value = 8018.1819
$ws.Cells.Item(122, 12).Value = 3900
$ws.Cells.Item(122, 13).Value = -5568.1819
$ws.Cells.Item(122, 14).Value = -8800
$ws.Cells.Item(126, 8).Value = 1527.2727
$ws.Cells.Item(126, 9).Value = 1200
$ws.Cells.Item(126, 10).Value = 1714.2858
$ws.Cells.Item(126, 11).Value = 3600
$ws.Cells.Item(126, 12).Value = 5142.857400000001
$ws.Cells.Item(126, 13).Value = -1130
$ws.Cells.Item(126, 14).Value = -10082.8574
$ws.Cells.Item(132, 8).Value = 26966.049
$ws.Cells.Item(132, 9).Value = 2574.682
$ws.Cells.Item(132, 10).Value = 55208.684
$ws.Cells.Item(132, 11).Value = 7724.045999999999
$ws.Cells.Item(132, 12).Value = 165626.052
$ws.Cells.Item(132, 13).Value = -5194.045999999999
$ws.Cells.Item(132, 14).Value = -170686.052

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(76, 8).Value = 20289.285
$ws.Cells.Item(76, 10).Value = 20289.285
$ws.Cells.Item(76, 12).Value = 20289.285
$ws.Cells.Item(76, 14).Value = -20965.285
$ws.Cells.Item(79, 8).Value = 20289.285
$ws.Cells.Item(79, 10).Value = 20289.285
$ws.Cells.Item(79, 12).Value = 20289.285
$ws.Cells.Item(79, 14).Value = -22629.285
$ws.Cells.Item(132, 8).Value = 573391.5
$ws.Cells.Item(132, 9).Value = 223856.44
$ws.Cells.Item(132, 10).Value = 835542.75
$ws.Cells.Item(132, 11).Value = 671569.3200000001
$ws.Cells.Item(132, 12).Value = 2506628.25
$ws.Cells.Item(132, 13).Value = -669039.3200000001
$ws.Cells.Item(132, 14).Value = -2511688.25
$ws.Cells.Item(136, 8).Value = 771647.5600000001
$ws.Cells.Item(136, 9).Value = 1430781.1
$ws.Cells.Item(136, 10).Value = 2658.3333
$ws.Cells.Item(136, 11).Value = 4292343.300000001
$ws.Cells.Item(136, 12).Value = 7974.999899999999
$ws.Cells.Item(136, 13).Value = -4289793.300000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4186.641
$ws.Cells.Item(132, 9).Value = 1418.3704
$ws.Cells.Item(132, 10).Value = 10415.25
$ws.Cells.Item(132, 11).Value = 4255.1112
$ws.Cells.Item(132, 12).Value = 31245.75
$ws.Cells.Item(132, 13).Value = -1725.1112
$ws.Cells.Item(132, 14).Value = -36305.75
$ws.Cells.Item(136, 8).Value = 590650
$ws.Cells.Item(136, 9).Value = 2710.3845
$ws.Cells.Item(136, 10).Value = 2501453.8
$ws.Cells.Item(136, 11).Value = 8131.1535
$ws.Cells.Item(136, 12).Value = 7504361.399999999
$ws.Cells.Item(136, 13).Value = -5581.1535
$ws.Cells.Item(136, 14).Value = -7509461.399999999
